$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price-column values that look like plain numbers (e.g. "218.72") get a
# leading apostrophe so Excel stores them as literal text instead of silently
# re-parsing them as a number (which would drop formatting like trailing zeros).

$ws.Range('D2').Value = '26.280.99'
$ws.Range('E2').Value = '  +0.66%  '
$ws.Range('D3').Value = '1.665.91'
$ws.Range('E3').Value = '  +0.73%  '
$ws.Range('E4').Value = '  +0.80%  '
$ws.Range('D5').Value = '''218.72'
$ws.Range('E5').Value = '  +0.56%  '
$ws.Range('D6').Value = '''0.5322'
$ws.Range('E6').Value = '  +1.51%  '
$ws.Range('D8').Value = '''0.2642'
$ws.Range('E8').Value = '  +1.54%  '
$ws.Range('D9').Value = '''0.06380'
$ws.Range('E9').Value = '  +0.68%  '
$ws.Range('D10').Value = '''20.54'
$ws.Range('E10').Value = '  +1.20%  '
$ws.Range('D11').Value = '''0.07825'
$ws.Range('E11').Value = '  +0.41%  '
$ws.Range('D12').Value = '''4.561'
$ws.Range('E12').Value = '  +1.41%  '
$ws.Range('D13').Value = '1.667.85'
$ws.Range('E13').Value = '  +0.98%  '
$ws.Range('D14').Value = '1.894.09'
$ws.Range('E14').Value = '  +0.69%  '
$ws.Range('D15').Value = '''0.5525'
$ws.Range('E15').Value = '  +1.17%  '
$ws.Range('D16').Value = '0.0₅8227'
$ws.Range('E16').Value = '  +0.55%  '
$ws.Range('D17').Value = '''65.70'
$ws.Range('E17').Value = '  +0.60%  '
$ws.Range('E18').Value = '  +0.84%  '
$ws.Range('D19').Value = '''4.687'
$ws.Range('E19').Value = '  +2.52%  '
$ws.Range('D20').Value = '''193.78'
$ws.Range('E20').Value = '  +1.45%  '
$ws.Range('E21').Value = '  +1.55%  '
$ws.Range('D22').Value = '''6.035'
$ws.Range('E22').Value = '  +0.43%  '
$ws.Range('E23').Value = '  +0.75%  '
$ws.Range('D24').Value = '''145.69'
$ws.Range('E24').Value = '  +2.51%  '
$ws.Range('D25').Value = '''0.1227'
$ws.Range('E25').Value = '  -0.74%  '
$ws.Range('E26').Value = '  -0.41%  '
$ws.Range('D27').Value = '''16.14'
$ws.Range('E27').Value = '  +0.31%  '
$ws.Range('D28').Value = '''1.486'
$ws.Range('E28').Value = '  +4.05%  '
$ws.Range('D29').Value = '''0.05896'
$ws.Range('E29').Value = '  +0.08%  '
$ws.Range('D30').Value = '''1.281'
$ws.Range('E30').Value = '  +0.30%  '
$ws.Range('D31').Value = '''3.609'
$ws.Range('E31').Value = '  +2.98%  '
$ws.Range('D32').Value = '''3.277'
$ws.Range('E32').Value = '  +1.08%  '
$ws.Range('D33').Value = '''1.609'
$ws.Range('E33').Value = '  +1.71%  '
$ws.Range('D34').Value = '''0.9636'
$ws.Range('E34').Value = '  +1.62%  '
$ws.Range('D35').Value = '''2.829'
$ws.Range('E35').Value = '  +1.78%  '
$ws.Range('E36').Value = '  +0.42%  '
$ws.Range('D37').Value = '''0.5804'
$ws.Range('E37').Value = '  +2.35%  '
$ws.Range('D38').Value = '''0.01608'
$ws.Range('E38').Value = '  -0.61%  '
$ws.Range('D39').Value = '''0.8662'
$ws.Range('E39').Value = '  +2.25%  '
$ws.Range('E40').Value = '  +0.50%  '
$ws.Range('D41').Value = '1.050.39'
$ws.Range('E41').Value = '  +2.43%  '
$ws.Range('D43').Value = '''104.07'
$ws.Range('E43').Value = '  +1.39%  '
$ws.Range('D44').Value = '1.804.83'
$ws.Range('E44').Value = '  +0.46%  '
$ws.Range('D45').Value = '''57.77'
$ws.Range('E45').Value = '  +1.26%  '
$ws.Range('D46').Value = '''1.014'
$ws.Range('E46').Value = '  +1.14%  '
$ws.Range('E47').Value = '  -5.81%  '
$ws.Range('D48').Value = '''8.076'
$ws.Range('E48').Value = '  +2.79%  '
$ws.Range('D49').Value = '''0.4380'
$ws.Range('E49').Value = '  +1.75%  '
$ws.Range('D50').Value = '''0.05163'
$ws.Range('E50').Value = '  -0.10%  '
$ws.Range('D51').Value = '''1.421'
$ws.Range('E51').Value = '  -3.20%  '
